# bug fixes on repost test 5
# Update column Q (Write Latency -> average) values for rows 3-23
# on the "tests results after queue" sheet.
#
# The cells hold numeric-looking values stored as text, so each new
# value is entered with a leading apostrophe to keep it as text
# (matching the original inline-string / text cell type) rather than
# letting Excel auto-convert it to a floating point number, which
# would silently drop significant trailing zeros (e.g. "8291.10").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = "8291.10"
    4  = "1911.77"
    5  = "15352.90"
    6  = "13191.98"
    7  = "1977.06"
    8  = "320754.56"
    9  = "8158.71"
    10 = "11553.72"
    11 = "5831.83"
    12 = "29.07"
    13 = "129147.26"
    14 = "25.49"
    15 = "2012.19"
    16 = "896051.34"
    17 = "16425.04"
    18 = "12891.88"
    19 = "105534.63"
    20 = "18.81"
    21 = "1916.08"
    22 = "130489.14"
    23 = "15297.58"
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = "'" + $updates[$row]
}
